# Weekly price update: a new price observation is inserted as the first
# dated record for this series (row 98), pushing the existing historical
# rows (old 98..141) down by one (new 99..142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 98 - shifts rows 98:141 down to 99:142
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly observation
$ws.Range("A98").Value = 1
$ws.Range("B98").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value = "Arica y Parinacota"
$ws.Range("D98").Value = 44846
$ws.Range("E98").Value = 15
$ws.Range("F98").Value = 100112036
$ws.Range("G98").Value = "Caigua"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 130
$ws.Range("K98").Value = 5000
$ws.Range("L98").Value = 6000
$ws.Range("M98").Value = 5500
$ws.Range("N98").Value = "$/caja 20 kilos"
$ws.Range("O98").Value = "Región de Arica y Parinacota"
$ws.Range("P98").Value = 275
$ws.Range("Q98").Value = 20
$ws.Range("R98").Value = "Hortaliza"
